$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new worksheet "2022-Q3" right after "总计" (so the tab
#    order becomes 总计, 2022-Q3, 2022-Q2, 2021-Q4).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $total)
$q3.Name = "2022-Q3"

# Style reference cells living on "总计" that already carry the bold /
# bordered / centered "header & index-column" look (style index 2 in the
# original workbook) so we can clone that exact formatting via copy/paste.
$headerStyleSrc = $total.Range("B1")
$indexStyleSrc = $total.Range("A2")

# --- header row ------------------------------------------------------------
$q3Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $q3Headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $q3Headers[$i]
}
$headerStyleSrc.Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

# --- data rows ---------------------------------------------------------------
# index | code | name | size | total position | position pct | held value | rank
$q3Rows = @(
    @(0, "009490", "泰康科技创新一年定期开放混合", "2.44", "79.62", "5.08", "0.1240", 4),
    @(1, "002935", "泰康恒泰回报灵活配置混合C", "2.53", "22.01", "1.50", "0.0380", 7),
    @(2, "002934", "泰康恒泰回报灵活配置混合A", "1.12", "22.01", "1.50", "0.0168", 7)
)

$r = 2
foreach ($row in $q3Rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = "'" + $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = "'" + $row[3]
    $q3.Cells.Item($r, 5).Value = "'" + $row[4]
    $q3.Cells.Item($r, 6).Value = "'" + $row[5]
    $q3.Cells.Item($r, 7).Value = "'" + $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}
$indexStyleSrc.Copy()
$q3.Range("A2:A4").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. "总计" gains a new row 2 for the 2022-Q3 totals; the previous rows 2/3
#    (2022-Q2 / 2021-Q4) shift down to rows 3/4 and their running index in
#    column A is bumped by one.
# ---------------------------------------------------------------------------
$total.Rows.Item(2).Insert()
# The freshly inserted row inherits stray formatting from the row above it
# (noticeable on B2:D2) - strip that back to the workbook default before
# writing the real 2022-Q3 values.
$total.Range("B2:D2").ClearFormats()
# Range refs are bound to a fixed A1 address, not live-tracking like real
# Excel: grab the "style 2" donor cell AFTER the shift, from what is now
# row 3 (the old, still-styled row 2).
$totalIndexStyleSrc = $total.Range("A3")

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 3
$total.Cells.Item(2, 4).Value = 0.18
$totalIndexStyleSrc.Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2

$excel.CutCopyMode = 0
